$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (duplicate header row), shifting rows 3-10 up to become rows 2-9.
$ws.Rows.Item(2).Delete()

# Clear the now-orphaned "Unnamed: 0" header text in A1.
$ws.Range("A1").Value = ""

# Remove the bold/bordered style that was applied to the header row.
$ws.Range("A1:P1").Style = "Normal"
